$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(559, 511, 424, 391, 327, 280, 206, 92, 0)
    3 = @(61, 49, 37, 34, 30, 26, 20, 13, 0)
    4 = @(165, 136, 104, 95, 75, 60, 45, 24, 0)
    5 = @(74, 51, 38, 32, 29, 21, 17, 11, 0)
    6 = @(19, 10, 9, 9, 8, 7, 4, 3, 0)
    7 = @(14, 2, 1, 1, 1, 1, 1, 1, 0)
    8 = @(11, 5, 2, 2, 1, 1, 1, 0, 0)
    9 = @(34, 109, 232, 276, 361, 409, 468, 578, 682)
    10 = @(7, 36, 50, 56, 59, 67, 71, 77, 88)
    11 = @(1, 12, 12, 12, 13, 14, 15, 16, 20)
    12 = @(11, 25, 33, 34, 37, 55, 91, 123, 148)
    13 = @(4, 14, 18, 18, 19, 19, 21, 22, 22)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
